$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = -5.873799999999997
$ws.Range("D4").Value  = -7.846900000000004
$ws.Range("E6").Value  = 12.46940000000001
$ws.Range("D7").Value  = -8.575399999999995
$ws.Range("E7").Value  = 11.7872
$ws.Range("D8").Value  = -8.694799999999994
$ws.Range("E8").Value  = 12.2753
$ws.Range("C11").Value = -13.32430000000001
$ws.Range("C12").Value = -14.76230000000002
$ws.Range("D12").Value = -8.153000000000006
$ws.Range("D14").Value = -8.564099999999998
$ws.Range("C15").Value = -11.9718
$ws.Range("E19").Value = 12.5675
$ws.Range("E21").Value = 12.63279999999999
$ws.Range("D22").Value = -8.342999999999996
$ws.Range("E24").Value = 12.77189999999999
$ws.Range("E25").Value = 13.16670000000001
